# edit.ps1 - applies the "Session 6 -> Session 7 / added material" edit
# described by the target diff to the RideShare/PHP lecture deck.
#
# Strategy: drive the PowerPoint object model exactly the way a human editor
# would - via TextRange.Characters(start, length).Text = "...". Re-assigning
# the text of a character range that spans more than one existing run causes
# PowerPoint to collapse that range into a single run (merging), while
# re-assigning a range that exactly matches one run but with different text
# causes it to retype that run (and, when the new text itself needs to be
# partially reformatted/retagged, PowerPoint splits it back into runs). This
# mirrors exactly the run-splitting / run-merging shown in the diff.
#
# NOTE: all mutations are written inline (no helper functions) - passing the
# live COM TextRange objects through PowerShell function parameters in this
# interpreter does not preserve the underlying mutation, so every
# Characters(...).Text = "..." assignment is performed directly at the call
# site.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1: title slide - "Session 6: PHP" -> "Session 7: PHP"
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$titleShape = $s1.Shapes.Item(3)
$titleRange = $titleShape.TextFrame.TextRange

# "6" -> "7" (single character run)
$f = $titleRange.Find("6")
$rng = $titleRange.Characters($f.Start, $f.Length)
$rng.Text = "7"

# Re-split ": PHP" into ": " and "PHP" so PHP becomes its own run
$f = $titleRange.Find("PHP")
$rng = $titleRange.Characters($f.Start, $f.Length)
$rng.Text = "PHP"

# ---------------------------------------------------------------------
# Slide 10: RideShare Exercise
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)

# Title: merge " " + "Exercise" runs into " Exercise"
$titleShape10 = $s10.Shapes.Item(1)
$titleRange10 = $titleShape10.TextFrame.TextRange
$f = $titleRange10.Find(" Exercise")
$rng = $titleRange10.Characters($f.Start, $f.Length)
$rng.Text = " Exercise"

# Body placeholder
$bodyShape10 = $s10.Shapes.Item(2)
$bodyRange10 = $bodyShape10.TextFrame.TextRange

# Para 1: merge the two runs into one sentence
$f = $bodyRange10.Find("Design a database to match drivers with passengers ")
$start = $f.Start
$f2 = $bodyRange10.Find("for ride sharing on long car trips:")
$len = ($f2.Start + $f2.Length) - $start
$rng = $bodyRange10.Characters($start, $len)
$rng.Text = "Design a database to match drivers with passengers for ride sharing on long car trips:"

# Para 3: "Passengers come looking..." + "and can make reservations" -> one run
$f = $bodyRange10.Find("Passengers come looking for rides: they want to know about available rides ")
$start = $f.Start
$f2 = $bodyRange10.Find("and can make reservations")
$len = ($f2.Start + $f2.Length) - $start
$rng = $bodyRange10.Characters($start, $len)
$rng.Text = "Passengers come looking for rides: they want to know about available rides and can make reservations"

# Para 4: "These things happen in no particular " + "order" -> one run
$f = $bodyRange10.Find("These things happen in no particular ")
$start = $f.Start
$f2 = $bodyRange10.Find("order")
$len = ($f2.Start + $f2.Length) - $start
$rng = $bodyRange10.Characters($start, $len)
$rng.Text = "These things happen in no particular order"

# Para 6 (last): "Build a " + "web application to accomplish the above" -> one run
$f = $bodyRange10.Find("Build a ")
$start = $f.Start
$f2 = $bodyRange10.Find("web application to accomplish the above")
$len = ($f2.Start + $f2.Length) - $start
$rng = $bodyRange10.Characters($start, $len)
$rng.Text = "Build a web application to accomplish the above"

# ---------------------------------------------------------------------
# Slide 11: RideShare Exercise: Tasks
# ---------------------------------------------------------------------
$s11 = $p.Slides.Item(11)

# Title: merge " Exercise: " + "Tasks" runs into " Exercise: Tasks"
$titleShape11 = $s11.Shapes.Item(1)
$titleRange11 = $titleShape11.TextFrame.TextRange
$f = $titleRange11.Find(" Exercise: ")
$start = $f.Start
$f2 = $titleRange11.Find("Tasks")
$len = ($f2.Start + $f2.Length) - $start
$rng = $titleRange11.Characters($start, $len)
$rng.Text = " Exercise: Tasks"

# Body: "What happens when a driver comes to find out who " + "the " + "passengers are?" -> one run
$bodyShape11 = $s11.Shapes.Item(2)
$bodyRange11 = $bodyShape11.TextFrame.TextRange
$f = $bodyRange11.Find("What happens when a driver comes to find out who ")
$start = $f.Start
$f2 = $bodyRange11.Find("passengers are?")
$len = ($f2.Start + $f2.Length) - $start
$rng = $bodyRange11.Characters($start, $len)
$rng.Text = "What happens when a driver comes to find out who the passengers are?"

# ---------------------------------------------------------------------
# Slide 4: Websites that are really databases
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$bodyShape4 = $s4.Shapes.Item(2)
$bodyRange4 = $bodyShape4.TextFrame.TextRange
$f = $bodyRange4.Find("Web pages are dynamically constructed ")
$start = $f.Start
$f2 = $bodyRange4.Find("results of database queries")
$len = ($f2.Start + $f2.Length) - $start
$rng = $bodyRange4.Characters($start, $len)
$rng.Text = "Web pages are dynamically constructed from results of database queries"

# ---------------------------------------------------------------------
# Slide 7: What is PHP?
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$bodyShape7 = $s7.Shapes.Item(2)
$bodyRange7 = $bodyShape7.TextFrame.TextRange

# "PHP is a " + "server-side " + "scripting language" -> one run
$f = $bodyRange7.Find("PHP is a ")
$start = $f.Start
$f2 = $bodyRange7.Find("scripting language")
$len = ($f2.Start + $f2.Length) - $start
$rng = $bodyRange7.Characters($start, $len)
$rng.Text = "PHP is a server-side scripting language"

# "More specifically, r" + "uns " + "inside the web " + "server" -> one run
$f = $bodyRange7.Find("More specifically, r")
$start = $f.Start
$f2 = $bodyRange7.Find("inside the web ")
# locate the final "server" occurrence after f2
$f3 = $bodyRange7.Find("server", $f2.Start)
$len = ($f3.Start + $f3.Length) - $start
$rng = $bodyRange7.Characters($start, $len)
$rng.Text = "More specifically, runs inside the web server"

# ---------------------------------------------------------------------
# Slide 8: PHP Scripts
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$bodyShape8 = $s8.Shapes.Item(2)
$bodyRange8 = $bodyShape8.TextFrame.TextRange
$f = $bodyRange8.Find("Are ")
$start = $f.Start
$f2 = $bodyRange8.Find("just like normal HTML pages")
$len = ($f2.Start + $f2.Length) - $start
$rng = $bodyRange8.Characters($start, $len)
$rng.Text = "Are just like normal HTML pages"
